$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row strings: "_old" -> "_FV2404" and "_new" -> "_FV2410"
$headers = @(
    "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
    "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
    "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn the used range into an Excel Table ("Table1") with an AutoFilter
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U76"), 0, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# 3) Freeze the header row (split below row 1) and keep the selection in the
#    lower pane, matching the frozen-pane view added to the sheet.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
